$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch C1 / C4 so they get an explicit (but unstyled) cell record, matching
# the "empty" placeholder cells that appear in column C once its other
# cells lose their inherited text-format style below.
$ws.Range("C1").Style = $ws.Range("A1").Style
$ws.Range("C4").Style = $ws.Range("A4").Style

# Reword the trait/quality notes in column C (rows 2, 3, 5) and drop their
# explicit "text" number-format style so they fall back to the default style.
$ws.Range("C2").Value = "Удачливый, старый"
$ws.Range("C2").Style = $ws.Range("A2").Style

$ws.Range("C3").Value = "Юркий, слабый"
$ws.Range("C3").Style = $ws.Range("A3").Style

$ws.Range("C5").Value = "бессмертный"
$ws.Range("C5").Style = $ws.Range("A5").Style

# Fix casing: "Кольцо" -> "кольцо" (row 7, "несет" relation)
$ws.Range("C7").Value = "кольцо"

# Insert a new row before the old row 10, shifting it down to row 11, and
# fill in the new relation: Фродо убегает от Назгул
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "Фродо"
$ws.Range("B10").Value = "убегает от"
$ws.Range("C10").Value = "Назгул"
$ws.Range("C10").Style = $ws.Range("A10").Style

[void]$ws.Range("D13").Select()
